# Applies the "Removed unused methods Updated TestData" change to TestData.xlsx
# Target sheet is the "DataSet" tab (1st tab) which holds the MSP Create
# customer API test rows (Kuwait = row2/3, UAE = row4/5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")
$ws.Activate()

# --- Row 2 (Kuwait / RunAPIRequest) -----------------------------------
# RunFlag Y -> N
$ws.Range("C2").Value = "N"

# Headers: refreshed Authorization token (new timestamp/signature)
$ws.Range("K2").Value = 'Authorization:osnAuth osnauth_x_application_id=6,  osnauth_x_source_id=14, osnauth_x_timestamp=1547012791, osnauth_x_signature=M2IxNTM4ZTU3NjRhNTk3OWRiYTA1YjI4Zjc3NDkxOTVlNmEwNjgxZWExNTc2M2Q4ZTE5MDA4ZDkzYjVhY2RhYg=='

# Body: MobileNumber updated from 96558880449075 to 96558880449081
$ws.Range("M2").Value = "{""UserId"":""00212029-ba97-468f-b670-b21eb2a93a8e"",`n""EmailAddress"":""info@osn.com"",`n""MobileNumber"":""96558880449081"",`n""Packages"":[  `n   3507`n],`n""BirthDate"":null,`n""Address"":null,`n""Country"":null,`n""City"":null,`n""Gender"":null,`n""CustomerUsernameID"":null,`n""Password"":""413703"",`n""Name"":null,`n""Title"":null,`n""LanguagePreference"":null,`n""Email2"":null,`n""Mobile2"":null,`n""Extra"":{  `n   ""MCC"":""965"",`n   ""MNC"":""01"",`n   ""Prod"":""01""`n},`n""CreatedDate"":""2018-12-12T13:00:54.4150898Z"",`n""ExpiryDate"":""2023-12-12T13:00:54.415093Z""}"

# ApiResponse: fresh userId + productExpiry timestamp
$ws.Range("N2").Value = '{"userId":"18823fc7-5cf8-4d0d-94ce-0c542f34c748","customerType":"OTT Telco","customerStatus":"OTT Active","accountCollection":[{"accountType":"OTT Msp","agreementCollection":[{"agreementType":877,"productCollection":[{"productId":3599,"productStatus":430,"productExpiry":"2069-01-09T05:13:50.420393+00:00"}]}]}],"responseCode":0,"returnId":0,"messageResponse":{"messageCode":0,"exceptionCode":0,"userMessages":null}}'

# --- Row 3 (Kuwait / ValidateApiResponse) ------------------------------
# RunFlag Y -> N
$ws.Range("C3").Value = "N"

# --- Row 4 (UAE / RunAPIRequest) ---------------------------------------
# Headers: same refreshed Authorization token as row 2
$ws.Range("K4").Value = 'Authorization:osnAuth osnauth_x_application_id=6,  osnauth_x_source_id=14, osnauth_x_timestamp=1547012791, osnauth_x_signature=M2IxNTM4ZTU3NjRhNTk3OWRiYTA1YjI4Zjc3NDkxOTVlNmEwNjgxZWExNTc2M2Q4ZTE5MDA4ZDkzYjVhY2RhYg=='

# Body: MobileNumber updated from 97110001008 to 971569143419
$ws.Range("M4").Value = '{"MobileNumber" : "971569143419", "EmailAddress" : "autoexection@osn.com", "Packages" : [3507], "Password" : "413703","extra": { "MCC": "971","MNC": "01","Prod": "01"}}'

# ApiResponse: fresh userId + productExpiry timestamp
$ws.Range("N4").Value = '{"userId":"90900888-eb50-45b9-978b-5747c6dcf71b","customerType":"OTT Telco","customerStatus":"OTT Active","accountCollection":[{"accountType":"OTT Msp","agreementCollection":[{"agreementType":877,"productCollection":[{"productId":3599,"productStatus":430,"productExpiry":"2069-01-09T06:04:52.238623+00:00"}]}]}],"responseCode":0,"returnId":0,"messageResponse":{"messageCode":0,"exceptionCode":0,"userMessages":null}}'

# --- Sheet view / selection --------------------------------------------
# Scroll back to top-left and move the active selection to K4.
$ws.Range("K4").Select()
